# The weekly refresh re-sorts the daily "Camote" price records (rows 2-37,
# columns A-R) by date. Every row's content moves as a whole to a new row;
# no cell is edited in place. Build the before->after row-destination map,
# snapshot all source rows first (so we never read a value that has already
# been overwritten), then write every row to its destination in one pass.

$map = @{2=17; 3=13; 4=34; 5=3; 6=4; 7=30; 8=23; 9=29; 10=8; 11=9; 12=27; 13=21; 14=15; 15=37; 16=33; 17=5; 18=6; 19=25; 20=16; 21=36; 22=20; 23=18; 24=35; 25=24; 26=7; 27=22; 28=31; 29=28; 30=12; 31=2; 32=19; 33=14; 34=11; 35=10; 36=26; 37=32}

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$firstRow = 2
$lastRow = 37
$firstCol = 1
$lastCol = 18

# Snapshot every source row's values (A:R) before writing anything.
$snapshot = @{}
for ($r = $firstRow; $r -le $lastRow; $r++) {
    $rowVals = @()
    for ($c = $firstCol; $c -le $lastCol; $c++) {
        $rowVals += $ws.Cells.Item($r, $c).Value2
    }
    $snapshot[$r] = $rowVals
}

# Write each snapshot row out to its mapped destination row.
foreach ($srcRow in $snapshot.Keys) {
    $dstRow = $map[$srcRow]
    $rowVals = $snapshot[$srcRow]
    for ($c = $firstCol; $c -le $lastCol; $c++) {
        $ws.Cells.Item($dstRow, $c).Value = $rowVals[$c - 1]
    }
}
